# 9th Stab- Cosmetic Changes
# Insert two new "week" columns (Jun_17, Jun_15) right after column A,
# pushing the existing Jun_13 / Jun_10 columns two slots to the right,
# and backfill the two new columns with the "UN" placeholder used
# throughout the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing B:C columns to D:E by inserting two blank columns
# in front of them.
$ws.Range("B:C").Insert()

# The displaced "week" column (old C, width 8) keeps its width; the two
# freshly inserted columns pick up the same width so all three line up.
$ws.Columns("C:E").ColumnWidth = 7.14

# New header row: newest week first.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Figure out how many data rows exist (column A holds the analyst / row
# labels for every data row below the header).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}
